$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 772-773), pushing the
# existing data (old rows 772..870) down to become rows 774..872.
$ws.Rows("772:773").Insert()

# Populate the new row 772 (weekly update: new date, "Primera" quality)
$ws.Range("A772").Value = 9
$ws.Range("B772").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C772").Value = "Metropolitana"
$ws.Range("D772").Value = 45077
$ws.Range("E772").Value = 13
$ws.Range("F772").Value = 100112009
$ws.Range("G772").Value = "Acelga"
$ws.Range("H772").Value = "Sin especificar"
$ws.Range("I772").Value = "Primera"
$ws.Range("J772").Value = 61
$ws.Range("K772").Value = 12000
$ws.Range("L772").Value = 12000
$ws.Range("M772").Value = 12000
$ws.Range("N772").Value = "$/docena de atados"
$ws.Range("O772").Value = "Región Metropolitana"
$ws.Range("P772").Value = 4000
$ws.Range("Q772").Value = 3
$ws.Range("R772").Value = "Hortaliza"

# Populate the new row 773 (weekly update: new date, "Segunda" quality)
$ws.Range("A773").Value = 9
$ws.Range("B773").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C773").Value = "Metropolitana"
$ws.Range("D773").Value = 45077
$ws.Range("E773").Value = 13
$ws.Range("F773").Value = 100112009
$ws.Range("G773").Value = "Acelga"
$ws.Range("H773").Value = "Sin especificar"
$ws.Range("I773").Value = "Segunda"
$ws.Range("J773").Value = 43
$ws.Range("K773").Value = 9000
$ws.Range("L773").Value = 9000
$ws.Range("M773").Value = 9000
$ws.Range("N773").Value = "$/docena de atados"
$ws.Range("O773").Value = "Región Metropolitana"
$ws.Range("P773").Value = 3000
$ws.Range("Q773").Value = 3
$ws.Range("R773").Value = "Hortaliza"
